$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused rows 13-36 first
for ($row = 36; $row -ge 13; $row--) {
    $ws.Rows.Item($row).Delete()
}

# Write rows 3..12 first (in row order), then row 2 last,
# matching the shared-string allocation order observed in the target file.
$ws.Range("A3").Value  = "生日快乐哈哈哈1"
$ws.Range("A4").Value  = "生日快乐2"
$ws.Range("A5").Value  = "生日快乐哈哈哈3"
$ws.Range("A6").Value  = "生日快乐4"
$ws.Range("A7").Value  = "生日快乐哈哈哈5"
$ws.Range("A8").Value  = "生日快乐6"
$ws.Range("A9").Value  = "生日快乐哈哈哈7"
$ws.Range("A10").Value = "生日快乐8"
$ws.Range("A11").Value = "生日快乐9"
$ws.Range("A12").Value = "生日快乐哈哈哈10"
$ws.Range("A2").Value  = "生日快乐0"

$ws.Range("A6").Select()
